# Add new worksheet "ihlp5050ez01" after "ihlp5050ce01" (before "ihlp6767dz01"),
# carrying the transposed/new mosfet_data table for the IHLP-5050EZ-01 part,
# and make it the active sheet (matches commit: "add ihlp, transpose mosfet_data").

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("ihlp5050ce01")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "ihlp5050ez01"

# Header row (column headers reuse the workbook's existing shared strings).
$headers = @("Lout","DCR","Iheat","Isat","Rth","Pheat","ET100","K0","K1","Kf","Kb","family","x","y","z")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows 2-19 (18 rows), columns A:O.
# Column index 11 (L / "family") is text; everything else is numeric.
$data = @(
    @("0.1","5.0000000000000001E-4","55","118","21.62","1.85","0.17","5.79","1.6100000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.22","5.9999999999999995E-4","51","110","20.82","1.92","0.92","53.68","1.6100000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.33","8.9999999999999998E-4","42","80","23.12","1.73","0.95","32.5","1.6100000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.47","1.1000000000000001E-3","38","65","21.82","1.83","1.1100000000000001","29","4.8300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.56000000000000005","1.2999999999999999E-3","36","55","20.57","1.94","1","23.04","4.3E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.68","1.5E-3","34","54","19.989999999999998","2","1.55","42.08","4.3E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("0.82","2E-3","31","53","18.03","2.2200000000000002","2.06","57.33","3.8300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("1","2.0999999999999999E-3","29","50","19.63","2.04","1.96","41.82","3.8300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("1.5","3.3999999999999998E-3","23","48","19.27","2.08","2.88","60.16","3.0300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("2.2000000000000002","4.5999999999999999E-3","20","32","18.84","2.12","2.91","49.82","3.0300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("3.3","7.7000000000000002E-3","15","32","20.010000000000002","2","2.9","41.87","3.0300000000000001E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("4.7","1.2800000000000001E-2","12","27","18.809999999999999","2.13","4.49","65.069999999999993","2.3999999999999998E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("5.6","1.4E-2","11.5","22","18.72","2.14","5.19","55.68","2.1299999999999999E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("6.8","1.54E-2","11","21","18.600000000000001","2.15","5.63","61.71","2.1299999999999999E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("7.8","1.72E-2","10","18","20.149999999999999","1.98","5.65","54.15","2.1299999999999999E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("8.1999999999999993","1.89E-2","9.5","18","20.32","1.97","5.67","39.43","7.8200000000000006E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("10","2.1399999999999999E-2","9","16","20","2","5.88","38.770000000000003","7.8200000000000006E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5"),
    @("15","2.8000000000000001E-2","8.25","14.5","18.190000000000001","2.2000000000000002","5.69","34.19","6.9699999999999996E-3","1.1879999999999999","2.1179999999999999","ihlp5050ez01","12.9","13.58","5")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        if ($c -eq 11) {
            $cell.Value = $row[$c]
        } else {
            $cell.Value = [double]$row[$c]
        }
    }
}

# Column L (12th column) width, matching the source sheet's custom width.
$ws.Range("L1").ColumnWidth = 13.7109375

# Leave the same selection/active-cell state recorded for this sheet, and make
# it the active (tab-selected) sheet/window, mirroring the recorded diff.
$ws.Range("N24").Select()
$ws.Activate()
